$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '26.293.73'
Set-TextValue $ws.Range('E2') '  +0.36%  '
Set-TextValue $ws.Range('D3') '1.589.38'
Set-TextValue $ws.Range('E3') '  +0.44%  '
Set-TextValue $ws.Range('E4') '  -0.19%  '
Set-TextValue $ws.Range('D5') '211.89'
Set-TextValue $ws.Range('E5') '  +1.06%  '
Set-TextValue $ws.Range('E6') '  +0.93%  '
Set-TextValue $ws.Range('E7') '  -0.19%  '
Set-TextValue $ws.Range('E8') '  +0.06%  '
Set-TextValue $ws.Range('E9') '  -0.22%  '
Set-TextValue $ws.Range('D10') '19.36'
Set-TextValue $ws.Range('E10') '  -0.74%  '
Set-TextValue $ws.Range('D11') '0.0847'
Set-TextValue $ws.Range('E11') '  +0.07%  '
Set-TextValue $ws.Range('D12') '1.813.38'
Set-TextValue $ws.Range('B13') 'Polkadot'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D13') '4.04'
Set-TextValue $ws.Range('E13') '  +0.07%  '
Set-TextValue $ws.Range('B14') 'WrappedEther'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D14') '1.576.16'
Set-TextValue $ws.Range('E14') '  -0.78%  '
Set-TextValue $ws.Range('E15') '  +0.87%  '
Set-TextValue $ws.Range('D16') '64.40'
Set-TextValue $ws.Range('E16') '  -0.17%  '
Set-TextValue $ws.Range('D17') '26.309.38'
Set-TextValue $ws.Range('E17') '  +0.42%  '
Set-TextValue $ws.Range('D18') '0.0₃0731'
Set-TextValue $ws.Range('E18') '  -0.38%  '
Set-TextValue $ws.Range('D19') '7.51'
Set-TextValue $ws.Range('E19') '  +3.45%  '
Set-TextValue $ws.Range('D20') '213.02'
Set-TextValue $ws.Range('D22') '4.28'
Set-TextValue $ws.Range('E22') '  +0.74%  '
Set-TextValue $ws.Range('E23') '  +1.41%  '
Set-TextValue $ws.Range('E24') '  -2.59%  '
Set-TextValue $ws.Range('D25') '145.20'
Set-TextValue $ws.Range('E25') '  +0.23%  '
Set-TextValue $ws.Range('E26') '  -0.14%  '
Set-TextValue $ws.Range('E27') '  +0.13%  '
Set-TextValue $ws.Range('E28') '  -0.51%  '
Set-TextValue $ws.Range('D29') '15.18'
Set-TextValue $ws.Range('E29') '  -0.16%  '
Set-TextValue $ws.Range('E30') '  -0.73%  '
Set-TextValue $ws.Range('E31') '  +0.89%  '
Set-TextValue $ws.Range('E32') '  -0.20%  '
Set-TextValue $ws.Range('E33') '  +1.09%  '
Set-TextValue $ws.Range('D34') '1.340.08'
Set-TextValue $ws.Range('E34') '  +4.48%  '
Set-TextValue $ws.Range('D35') '2.44'
Set-TextValue $ws.Range('E36') '  -0.92%  '
Set-TextValue $ws.Range('E37') '  +0.09%  '
Set-TextValue $ws.Range('E38') '  -0.06%  '
Set-TextValue $ws.Range('D39') '1.06'
Set-TextValue $ws.Range('E39') '  -14.02%  '
Set-TextValue $ws.Range('D40') '0.817'
Set-TextValue $ws.Range('E40') '  +0.37%  '
Set-TextValue $ws.Range('E41') '  +3.75%  '
Set-TextValue $ws.Range('E42') '  -0.14%  '
Set-TextValue $ws.Range('E43') '  +0.41%  '
Set-TextValue $ws.Range('D44') '0.762'
Set-TextValue $ws.Range('E44') '  -0.56%  '
Set-TextValue $ws.Range('D45') '1.725.07'
Set-TextValue $ws.Range('E45') '  +0.34%  '
Set-TextValue $ws.Range('D46') '61.71'
Set-TextValue $ws.Range('E46') '  -0.89%  '
Set-TextValue $ws.Range('D47') '88.07'
Set-TextValue $ws.Range('E47') '  -0.86%  '
Set-TextValue $ws.Range('B48') 'RenderToken'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D48') '1.49'
Set-TextValue $ws.Range('E48') '  -4.02%  '
Set-TextValue $ws.Range('B49') 'Algorand'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D49') '0.0980'
Set-TextValue $ws.Range('E49') '  -2.59%  '
Set-TextValue $ws.Range('B50') 'Cronos'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D50') '0.0503'
Set-TextValue $ws.Range('E50') '  -0.64%  '
Set-TextValue $ws.Range('B51') 'USDD'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws.Range('D51') '0.999'
Set-TextValue $ws.Range('E51') '  -0.40%  '
